$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E9 was stored as a genuine number (20000); the edit changes it to a
# text value "20000" instead (matching the other rows in column E that
# already store the amount as text). Use a leading apostrophe so the
# engine treats the input as text rather than re-parsing it as a
# number, then reset the style back to Normal so no stray
# quote-prefix / text-number-format style gets attached to the cell.
$ws.Range("E9").Value = "'20000"
$ws.Range("E9").Style = "Normal"

# Append the new payment row (row 10).
$ws.Range("A10").Value = "'9801234567819235"
$ws.Range("A10").Style = "Normal"

$ws.Range("B10").Value = "BG8799BI"
$ws.Range("C10").Value = "Dinda"
$ws.Range("D10").Value = "01-08-2025 08:37"

# E10 stays a genuine number (10000), unlike E9.
$ws.Range("E10").Value = 10000

$ws.Range("F10").Value = "Bank Mandiri"
